$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string rich text edits (preserve approximate text) ---
$ws.Range("A8").Characters(21, 2).Text = "50"
$ws.Range("C9").Characters(47, 10).Text = "12/17/2023"
$ws.Range("C9").Characters(27, 9).Text = "12/11/2023"

# --- Cells that change type (number <-> text) ---
# Use PasteSpecial(xlPasteFormats) from a donor cell of the right existing
# style, then set value, so the canonical style index is reused exactly.
$donorNum = $ws.Range("C36")   # plain number style (15)
$donorTxt = $ws.Range("A14")   # plain text style (14)

$ws.Range("G15").Value = "'0"
$donorTxt.Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "'***.*"
$donorTxt.Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("C23").Value = "'0"
$donorTxt.Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("D23").Value = "'0"
$donorTxt.Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = "'***.*"
$donorTxt.Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("C26").Value = 1
$donorNum.Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("F26").Value = 1
$donorNum.Copy()
$ws.Range("F26").PasteSpecial(-4122)

$ws.Range("G26").Value = "'0"
$donorTxt.Copy()
$ws.Range("G26").PasteSpecial(-4122)

$ws.Range("H26").Value = "'***.*"
$donorTxt.Copy()
$ws.Range("H26").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$donorTxt.Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = "'0"
$donorTxt.Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "'***.*"
$donorTxt.Copy()
$ws.Range("E27").PasteSpecial(-4122)

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("L15").Value = -50
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 115
$ws.Range("J16").Value = 137
$ws.Range("K16").Value = -16.058394160583
$ws.Range("L16").Value = 43.75
$ws.Range("M16").Value = -8
$ws.Range("N16").Value = -80.341880341880
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 109
$ws.Range("J17").Value = 155
$ws.Range("K17").Value = -29.677419354838
$ws.Range("L17").Value = 6.862745098039
$ws.Range("M17").Value = 53.521126760563
$ws.Range("N17").Value = -57.915057915057
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -85.714285714285
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -38.095238095238
$ws.Range("I18").Value = 167
$ws.Range("J18").Value = 258
$ws.Range("K18").Value = -35.271317829457
$ws.Range("L18").Value = -25.112107623318
$ws.Range("M18").Value = -28.936170212766
$ws.Range("N18").Value = -84.762773722627
$ws.Range("C19").Value = 16
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -17.307692307692
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = 595
$ws.Range("K19").Value = 9.243697478991
$ws.Range("L19").Value = 49.082568807339
$ws.Range("M19").Value = 120.338983050847
$ws.Range("N19").Value = 81.058495821727
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 153
$ws.Range("J20").Value = 188
$ws.Range("K20").Value = -18.617021276595
$ws.Range("L20").Value = 51.485148514851
$ws.Range("M20").Value = 6.993006993006
$ws.Range("N20").Value = -82.534246575342
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -14.814814814814
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -28.571428571428
$ws.Range("I21").Value = 1200
$ws.Range("J21").Value = 1347
$ws.Range("K21").Value = -10.913140311804
$ws.Range("L21").Value = 25.786163522012
$ws.Range("M21").Value = 37.931034482758
$ws.Range("N21").Value = -62.394233782513
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 44.444444444444
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = 11.290322580645
$ws.Range("I24").Value = 875
$ws.Range("J24").Value = 978
$ws.Range("K24").Value = -10.531697341513
$ws.Range("L24").Value = 0.114416475972
$ws.Range("M24").Value = 48.556876061120
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 235
$ws.Range("J25").Value = 269
$ws.Range("K25").Value = -12.639405204461
$ws.Range("L25").Value = 4.910714285714
$ws.Range("M25").Value = 15.763546798029
$ws.Range("I26").Value = 9
$ws.Range("K26").Value = -40
$ws.Range("L26").Value = -50
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = 46.666666666666
